# Apply the "DB.xlsx" recipe-table corrections to the "Recipes" sheet:
#  - Re-number the ID column sequentially (11..19), closing the gap left
#    by the missing recipe #5 in the original 1,2,3,4,6,7,8,9,10 sequence
#  - Fix the Cuisine / Skill / Time columns, which were shifted one column
#    to the left (Cuisine ended up in the Skill column, Skill in Time, and
#    Time had spilled into the Cuisine column of the next logical field).
#    The correct values are recovered by rotating C<-E, D<-C, E<-D.
#  - Row 6 and Row 7 actually belong to different recipes: row 6 held the
#    "Chicken Curry" recipe and row 7 held "Ban Mian (Soup)" - they need to
#    swap places (Ban Mian becomes row 6, Chicken Curry becomes row 7).
#  - Two Instruction cells (F3 and the Chicken Curry Instruction cell) had
#    an explicit word-wrap style applied; that formatting is removed.
#  - Row 9's Ingredient cell stored a Python list repr; convert it to a
#    plain comma separated string.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Recipes")

function Fix-CuisineSkillTime($row) {
    $c = $ws.Range("C$row").Value2
    $d = $ws.Range("D$row").Value2
    $e = $ws.Range("E$row").Value2
    $ws.Range("C$row").Value2 = $e
    $ws.Range("D$row").Value2 = $c
    $ws.Range("E$row").Value2 = $d
}

# --- Simple rows: only the ID and the Cuisine/Skill/Time columns move.
#     The IDs are renumbered sequentially (row number + 9) to close the
#     gap left by the missing recipe #5: 1,2,3,4,6,7,8,9,10 -> 11..19 ---
$simpleRows = @(2, 3, 4, 5, 8, 9, 10)
foreach ($r in $simpleRows) {
    $ws.Range("A$r").Value2 = $r + 9
    Fix-CuisineSkillTime $r
}

# --- Rows 6 & 7: the two recipes swap rows entirely (B..J), then each
#     gets its Cuisine/Skill/Time columns corrected the same way ---
$row6Vals = $ws.Range("B6:J6").Value2
$row7Vals = $ws.Range("B7:J7").Value2
$ws.Range("B6:J6").Value2 = $row7Vals
$ws.Range("B7:J7").Value2 = $row6Vals

$ws.Range("A6").Value2 = 15
$ws.Range("A7").Value2 = 16
Fix-CuisineSkillTime 6
Fix-CuisineSkillTime 7

# Pasting the (longer/shorter) Instruction text can make the engine
# auto-fit the row height; restore the original explicit row heights.
$ws.Rows.Item(6).RowHeight = 14.25
$ws.Rows.Item(7).RowHeight = 13.5

# --- Remove the word-wrap formatting from F3 and the Instruction cell
#     that now sits in F6 (it used to be F6 before the swap above).
#     Copying the (plain) style from F2 clears the wrap-text formatting
#     entirely, rather than leaving a redundant "no-op" style behind. ---
$ws.Range("F3").Style = $ws.Range("F2").Style
$ws.Range("F6").Style = $ws.Range("F2").Style

# --- Row 9 Ingredient cell: turn the Python-list-style text into a plain
#     comma separated list ---
$ws.Range("G9").Value2 = "240 ml Chicken Stock, 0.5 Stalk Lemongrass (Halved), 0.75 Galangal (Sliced), 10 g Tom Yum Chilli Paste, 0.5 Kaffir Lime Leaf, 3 Oyster Mushrooms, 0.25 Thai Chilli Pepper (Halved), 2 Large Shrimp, 3 g Sugar, 4 ml Lime Juice, 10 g Fresh Cilantro Leaves"
